$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add value for C25 (time-of-day value, styled like the other C cells)
$ws.Range("C25").Value = 0.56944444444444442
$ws.Range("C25").NumberFormat = $ws.Range("C24").NumberFormat

# Update the active cell selection to C26 (matching the new data entry point)
$ws.Range("C26").Select()
